$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 50, shifting existing rows 50-54 down to 51-55
$ws.Rows.Item(50).Insert()

# Fill in the new row 50 with the data for the new entry
$ws.Range("A50").Value = 10
$ws.Range("B50").Value = "Vega Modelo de Temuco"
$ws.Range("C50").Value = "La Araucanía"
$ws.Range("D50").Value = 44769
$ws.Range("E50").Value = 9
$ws.Range("F50").Value = 100112010
$ws.Range("G50").Value = "Achicoria"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 400
$ws.Range("K50").Value = 10000
$ws.Range("L50").Value = 11000
$ws.Range("M50").Value = 10750
$ws.Range("N50").Value = "$/caja 18 unidades"
$ws.Range("O50").Value = "Región Metropolitana"
$ws.Range("P50").Value = 597
$ws.Range("Q50").Value = 18
$ws.Range("R50").Value = "Hortaliza"
